# Borrar puntos al final de frases en viñetas para versionas Full CV
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("supervision")

# Remove trailing periods from bullet-point sentences in column E (and the
# "S." -> "S" abbreviation fix inside E4).
# Order matters for how new strings are appended to the shared string table,
# so update E11:E22 before E2:E5.
$ws.Range("E11").Value = "Paula Andrea Betancourt Velandia  (2018 - 2019)"
$ws.Range("E12").Value = "Ana Sofía Gómez Castelblanco (2018 - 2019)"
$ws.Range("E13").Value = "Lina María García Hoyos  (2016 - 2017)"
$ws.Range("E14").Value = "Angie Liliana Pérez Rodríguez  (2016 - 2018)"
$ws.Range("E15").Value = "Lina María Morales Sánchez (2016 - 2017)"
$ws.Range("E16").Value = "Laura Milena Estupiñan Aldana  (2016 - 2017)"
$ws.Range("E17").Value = "Vanesa Díaz Güiza  (2016 - 2018)"
$ws.Range("E18").Value = "Cindy Paola Moncada Gómez (2016 - 2017)"
$ws.Range("E19").Value = "Haydn Ricardo Roldán Morales (2015 - 2016)"
$ws.Range("E20").Value = "Maria Alejandra Abello Mozo  (2018 - 2018)"
$ws.Range("E21").Value = "Natalia Elízabeth Moreno Buitrago (2017 - 2019)"
$ws.Range("E22").Value = "Juan Felipe Pérez Ariza (2017 - 2019)"

$ws.Range("E2").Value = "Milena Vásquez-Amézquita. Supervisión conjunta con  Alicia Salvador"
$ws.Range("E3").Value = "Francisco Javier Flores. Supervisión conjunta con Lisa Chiara Fellin"
$ws.Range("E4").Value = "Julia Sanz-Vidania. Supervisión conjunta con S Craig Roberts"
$ws.Range("E5").Value = "Adrián Acosta Guerrero. Supervisión conjunta con Milena Vásquez-Amézquita"

$ws.Range("E22").Select()
